$d = $word.ActiveDocument

$replacements = @(
    @{old="65÷3=21, 2"; new="21÷9=2, 3"},
    @{old="52÷4=13, 0"; new="30÷4=7, 2"},
    @{old="50÷7=7, 1"; new="92÷6=15, 2"},
    @{old="98÷4=24, 2"; new="47÷4=11, 3"},
    @{old="58÷2=29, 0"; new="34÷5=6, 4"},
    @{old="61÷2=30, 1"; new="44÷2=22, 0"},
    @{old="32÷5=6, 2"; new="18÷4=4, 2"},
    @{old="81÷3=27, 0"; new="53÷8=6, 5"},
    @{old="32÷2=16, 0"; new="18÷6=3, 0"},
    @{old="93÷3=31, 0"; new="52÷5=10, 2"},
    @{old="31÷2=15, 1"; new="30÷8=3, 6"},
    @{old="72÷4=18, 0"; new="26÷5=5, 1"},
    @{old="41÷3=13, 2"; new="22÷2=11, 0"},
    @{old="19÷4=4, 3"; new="17÷9=1, 8"},
    @{old="18÷9=2, 0"; new="88÷8=11, 0"},
    @{old="24÷4=6, 0"; new="52÷9=5, 7"},
    @{old="35÷6=5, 5"; new="35÷5=7, 0"},
    @{old="66÷8=8, 2"; new="28÷8=3, 4"},
    @{old="73÷5=14, 3"; new="10÷5=2, 0"},
    @{old="39÷3=13, 0"; new="59÷2=29, 1"},
    @{old="48÷7=6, 6"; new="55÷7=7, 6"},
    @{old="33÷3=11, 0"; new="96÷9=10, 6"},
    @{old="37÷9=4, 1"; new="46÷3=15, 1"},
    @{old="54÷2=27, 0"; new="66÷4=16, 2"},
    @{old="57÷3=19, 0"; new="34÷5=6, 4"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
